$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtTest10mm")
[void]$ws.Activate()

$ws.Range("J21").Value = "Vas_Pam insertion wrt arbitrary knee loc"

$ws.Range("J22").Value = "x"
$ws.Range("K22").Value = "y"
$ws.Range("L22").Value = "z"

$ws.Range("J23").Value = 0.02163
$ws.Range("K23").Value = -0.07164
$ws.Range("L23").Value = 0

$ws.Range("J25").Value = "Vas_Pam insertion wrt arbitrary knee loc"

$ws.Range("J26").Value = "x"
$ws.Range("K26").Value = "y"
$ws.Range("L26").Value = "z"

$ws.Range("J27").Value = 0.03639
$ws.Range("K27").Value = -0.06313
$ws.Range("L27").Value = 0

[void]$ws.Range("K27").Select()
